$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell D-column values in this sheet are plain text (e.g. "37.774.20", "59.00")
# using "." as both thousands separator and a pseudo-decimal point. Several of the
# updated values (e.g. "234.21", "59.00") look like ordinary numbers, so the target
# cell is pre-formatted as Text ("@") before assignment to keep them literal strings
# (preserving exact formatting such as trailing zeros) instead of being parsed as numbers.

$ws.Range('D2').Value = '37.786.62'
$ws.Range('E2').Value = '  +1.06%  '
$ws.Range('D3').Value = '2.084.04'
$ws.Range('E3').Value = '  +0.75%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '234.21'
$ws.Range('E5').Value = '  -0.36%  '
$ws.Range('E6').Value = '  -0.48%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '59.00'
$ws.Range('E7').Value = '  +3.25%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.391'
$ws.Range('E9').Value = '  -0.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0788'
$ws.Range('E10').Value = '  +1.80%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.105'
$ws.Range('E11').Value = '  +2.76%  '
$ws.Range('D12').Value = '2.390.24'
$ws.Range('E12').Value = '  +0.82%  '
$ws.Range('E13').Value = '  +2.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.23'
$ws.Range('E14').Value = '  +2.73%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.774'
$ws.Range('E15').Value = '  -0.15%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.29'
$ws.Range('E16').Value = '  +2.06%  '
$ws.Range('D17').Value = '2.079.52'
$ws.Range('D18').Value = '37.699.32'
$ws.Range('E18').Value = '  +1.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.17'
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('E20').Value = '  +2.59%  '
$ws.Range('E21').Value = '  +1.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '228.65'
$ws.Range('E22').Value = '  +0.81%  '
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.41'
$ws.Range('E24').Value = '  -0.85%  '
$ws.Range('E25').Value = '  -0.92%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '170.16'
$ws.Range('E26').Value = '  +1.85%  '
$ws.Range('E27').Value = '  +8.31%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.02'
$ws.Range('E28').Value = '  +1.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.42'
$ws.Range('E29').Value = '  +0.78%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.53'
$ws.Range('E30').Value = '  +2.27%  '
$ws.Range('E31').Value = '  +1.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.72'
$ws.Range('E32').Value = '  +3.87%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.73'
$ws.Range('E33').Value = '  +4.35%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0629'
$ws.Range('E34').Value = '  +2.11%  '
$ws.Range('E35').Value = '  +1.99%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.45'
$ws.Range('E36').Value = '  +3.78%  '
$ws.Range('E37').Value = '  +2.60%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.999'
$ws.Range('E38').Value = '  -0.31%  '
$ws.Range('E39').Value = '  -3.77%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0990'
$ws.Range('E40').Value = '  +3.53%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.95'
$ws.Range('E41').Value = '  +0.12%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '98.74'
$ws.Range('E42').Value = '  +1.17%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0215'
$ws.Range('E43').Value = '  +1.07%  '
$ws.Range('D44').Value = '1.457.94'
$ws.Range('E44').Value = '  -2.10%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.32'
$ws.Range('E45').Value = '  +4.18%  '
$ws.Range('E46').Value = '  +1.07%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '16.02'
$ws.Range('E47').Value = '  +6.25%  '
$ws.Range('E48').Value = '  +4.40%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.40'
$ws.Range('E49').Value = '  +2.96%  '
$ws.Range('E50').Value = '  +2.53%  '
$ws.Range('D51').Value = '2.274.90'
$ws.Range('E51').Value = '  +0.72%  '
